# Update workbook/sheet metadata: rename sheet to reflect the new export
# timestamp (20240806-075647 -> 20240807-074757).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = "IClientBalance-20240807-074757-"

# Column G holds the balance date for every data row (2..274). The whole
# export was regenerated a day later, so every date moves from 2024-08-06
# (serial 45510) to 2024-08-07 (serial 45511).
for ($r = 2; $r -le 274; $r++) {
    $ws.Cells.Item($r, 7).Value = 45511
}

# A handful of rows also had their balance (column E) and total (column H)
# amounts corrected in the refreshed export. Apply those specific updates.
$ws.Cells.Item(17, 5).Value = 708.26
$ws.Cells.Item(17, 8).Value = 708.26

$ws.Cells.Item(52, 5).Value = 999.9
$ws.Cells.Item(52, 8).Value = 999.9

$ws.Cells.Item(60, 5).Value = 16821.66
$ws.Cells.Item(60, 8).Value = 16821.66

$ws.Cells.Item(103, 5).Value = 904.13
$ws.Cells.Item(103, 8).Value = 904.13

$ws.Cells.Item(245, 5).Value = 8462.0400000000009
$ws.Cells.Item(245, 8).Value = 8462.0400000000009
